$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "D:\Users\Chickens\Documents\EPCC\SynthSys\code_projects\synbio-toolkit\src\test\resources\ed\biordm\sbol\synbio\handler\"

$ws.Range("B3").Value = $prefix + "NC_035470.gbk"
$ws.Range("B5").Value = $prefix + "NC_014139.gbk"
$ws.Range("B2").Value = $prefix + "NC_001499.gbk"

[void]$ws.Range("D16").Select()
